$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cm flag for cm001 (row 2) from FALSE to TRUE
$ws.Range("C2").Value = $true

# Move the active selection to C3 as reflected in the saved view state
$ws.Range("C3").Select()
